# Fruta / hortaliza, semanal
# Insert a new weekly record at row 154, pushing the existing rows 154..229
# down to 155..230, and populate the new row 154 with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 154 (shifts rows 154-229 down to 155-230)
$ws.Rows.Item(154).Insert()

# Populate the newly inserted row 154 with the new record's values
$ws.Cells.Item(154, 1).Value = 7
$ws.Cells.Item(154, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(154, 3).Value = "Ñuble"
$ws.Cells.Item(154, 4).Value = 44609
$ws.Cells.Item(154, 5).Value = 16
$ws.Cells.Item(154, 6).Value = 100112008
$ws.Cells.Item(154, 7).Value = "Coliflor"
$ws.Cells.Item(154, 8).Value = "Sin especificar"
$ws.Cells.Item(154, 9).Value = "Primera"
$ws.Cells.Item(154, 10).Value = 200
$ws.Cells.Item(154, 11).Value = 800
$ws.Cells.Item(154, 12).Value = 850
$ws.Cells.Item(154, 13).Value = 825
$ws.Cells.Item(154, 14).Value = '$/unidad'
$ws.Cells.Item(154, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(154, 16).Value = 825
$ws.Cells.Item(154, 17).Value = 1
$ws.Cells.Item(154, 18).Value = "Hortaliza"
